$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 103
$ws.Range("C4").Value = "babby"

$ws.Range("C4").Select()
